$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- sharedStrings.xml changes realized as direct cell edits ---
# 1) Delete si "سرمايه گذاري صبا تامين" (was referenced by row 675) and
#    insert new si "دولت جمهوري اسلامي ايران" before "اختصاصی بازارگردانی ملت" (row 1075).
# Net effect on data rows 675-1074: every (Nationalcode, Name) pair shifts up by one row,
# and a brand-new entry ("دولت جمهوري اسلامي ايران" bank) is appended at row 1074.

# --- Column B (Nationalcode) and C (Name) updates for rows 675-1074 ---
$ws.Range("B675").Value = 10700144582
$ws.Range("C675").Value = 'پتروشيمي گلستان'
$ws.Range("B676").Value = 10220016318
$ws.Range("C676").Value = 'سيمان‌اروميه‌'
$ws.Range("B677").Value = 10460051496
$ws.Range("C677").Value = 'پارس‌سويچ‌'
$ws.Range("B678").Value = 10102239945
$ws.Range("C678").Value = 'صنايع بهداشتي ساينا'
$ws.Range("B679").Value = 10101212701
$ws.Range("C679").Value = 'همكاران سيستم'
$ws.Range("B680").Value = 10101754741
$ws.Range("C680").Value = 'سازه‌ پويش‌'
$ws.Range("B681").Value = 10260328876
$ws.Range("C681").Value = 'تابان نيرو سپاهان'
$ws.Range("B682").Value = 10100387143
$ws.Range("C682").Value = 'ايران‌ تاير'
$ws.Range("B683").Value = 10100580397
$ws.Range("C683").Value = 'تامين‌ ماسه‌ ريخته‌گري‌'
$ws.Range("B684").Value = 10101044869
$ws.Range("C684").Value = 'سرمايه گذاري تامين اجتماعي'
$ws.Range("B685").Value = 14005751499
$ws.Range("C685").Value = 'پيشگامان فن آوري و دانش آراميس'
$ws.Range("B686").Value = 10320839651
$ws.Range("C686").Value = 'بيمه تعاون'
$ws.Range("B687").Value = 10100478433
$ws.Range("C687").Value = 'تايدواترخاورميانه‌'
$ws.Range("B688").Value = 10200142417
$ws.Range("C688").Value = 'سرمايه‌گذاري‌توسعه‌آذربايجان‌'
$ws.Range("B689").Value = 10260357602
$ws.Range("C689").Value = 'ذغال‌سنگ‌ نگين‌ ط‌بس‌'
$ws.Range("B690").Value = 10260377008
$ws.Range("C690").Value = 'فرآوري زغال سنگ پروده طبس'
$ws.Range("B691").Value = 10101907666
$ws.Range("C691").Value = 'گروه سرمايه گذاري تدبير'
$ws.Range("B692").Value = 10260439863
$ws.Range("C692").Value = 'سرمايه گذاري توسعه توكا'
$ws.Range("B693").Value = 10380266337
$ws.Range("C693").Value = 'توسعه‌شهري‌توس‌گستر'
$ws.Range("B694").Value = 10100254685
$ws.Range("C694").Value = 'داروسازي‌ تهران‌ دارو'
$ws.Range("B695").Value = 10380466986
$ws.Range("C695").Value = 'گروه كارخانجات صنعتي تبرك'
$ws.Range("B696").Value = 10103251617
$ws.Range("C696").Value = 'صنايع پتروشيمي تخت جمشيد'
$ws.Range("B697").Value = 10260289240
$ws.Range("C697").Value = 'تكادو'
$ws.Range("B698").Value = 10101477062
$ws.Range("C698").Value = 'كنترل ‌خوردگي ‌تكين ‌كوي'
$ws.Range("B699").Value = 10000000004
$ws.Range("C699").Value = 'توليدكنندگان بورس كالاي ايران'
$ws.Range("B700").Value = 10861658680
$ws.Range("C700").Value = 'تكنوتار'
$ws.Range("B701").Value = 10260294158
$ws.Range("C701").Value = 'توليدي‌ كاشي‌ تكسرام‌'
$ws.Range("B702").Value = 10101021508
$ws.Range("C702").Value = 'دامداري تليسه نمونه'
$ws.Range("B703").Value = 10320741707
$ws.Range("C703").Value = 'تامين سرمايه تمدن'
$ws.Range("B704").Value = 10102501905
$ws.Range("C704").Value = 'سرمايه‌گذاري‌توسعه‌ملي‌'
$ws.Range("B705").Value = 10320585987
$ws.Range("C705").Value = 'بين المللي توسعه ص. معادن غدير'
$ws.Range("B706").Value = 10101136712
$ws.Range("C706").Value = 'توليدمحورخودرو'
$ws.Range("B707").Value = 10320453421
$ws.Range("C707").Value = 'تامين سرمايه بانك ملت'
$ws.Range("B708").Value = 10320826260
$ws.Range("C708").Value = 'تجلي توسعه معادن و فلزات'
$ws.Range("B709").Value = 10101725603
$ws.Range("C709").Value = 'توليدمواداوليه‌داروپخش‌'
$ws.Range("B710").Value = 14005725237
$ws.Range("C710").Value = 'بيمه تجارت نو'
$ws.Range("B711").Value = 10103584900
$ws.Range("C711").Value = 'تامين سرمايه نوين'
$ws.Range("B712").Value = 14005652664
$ws.Range("C712").Value = 'سرمايه گذاري اقتصاد شهر طوبي'
$ws.Range("B713").Value = 10260272518
$ws.Range("C713").Value = 'سرمايه‌گذاري‌توكافولاد(هلدينگ‌'
$ws.Range("B714").Value = 10260041071
$ws.Range("C714").Value = 'توسعه و عمران شهرستان نائين'
$ws.Range("B715").Value = 10101912330
$ws.Range("C715").Value = 'توسعه و عمران اميد'
$ws.Range("B716").Value = 10100595801
$ws.Range("C716").Value = 'تولي‌پرس‌'
$ws.Range("B717").Value = 10260450980
$ws.Range("C717").Value = 'توكاريل'
$ws.Range("B718").Value = 10101702561
$ws.Range("C718").Value = 'گروه س توسعه صنعتي ايران'
$ws.Range("B719").Value = 10103637809
$ws.Range("C719").Value = 'توسعه سامانه ي نرم افزاري نگين'
$ws.Range("B720").Value = 10840053064
$ws.Range("C720").Value = 'زغال سنگ پروده طبس'
$ws.Range("B721").Value = 10860726096
$ws.Range("C721").Value = 'پتروشيمي تندگويان'
$ws.Range("B722").Value = 10260322490
$ws.Range("C722").Value = 'مهندسي مرآت پولاد'
$ws.Range("B723").Value = 10260492543
$ws.Range("C723").Value = 'توكا رنگ فولاد سپاهان'
$ws.Range("B724").Value = 10200044047
$ws.Range("C724").Value = 'تراكتورسازي‌ايران‌'
$ws.Range("B725").Value = 10861521552
$ws.Range("C725").Value = 'توريستي ورفاهي آبادگران كيش'
$ws.Range("B726").Value = 10100433825
$ws.Range("C726").Value = 'ايران‌ ترانسفو'
$ws.Range("B727").Value = 10103525426
$ws.Range("C727").Value = 'تامين سرمايه امين'
$ws.Range("B728").Value = 10101117032
$ws.Range("C728").Value = 'توليد سموم‌ علف‌ كش‌'
$ws.Range("B729").Value = 10100539090
$ws.Range("C729").Value = 'توسعه‌ صنايع‌ بهشهر(هلدينگ‌'
$ws.Range("B730").Value = 10101767219
$ws.Range("C730").Value = 'اعتباري توسعه'
$ws.Range("B731").Value = 10100782085
$ws.Range("C731").Value = 'سرمايه گذاري پارس‌ توشه‌'
$ws.Range("B732").Value = 10100993098
$ws.Range("C732").Value = 'كارخانجات‌توليدي‌شيشه‌رازي‌'
$ws.Range("B733").Value = 10380019602
$ws.Range("C733").Value = 'قند تربت حيدريه'
$ws.Range("B734").Value = 10460096368
$ws.Range("C734").Value = 'ترانسفورماتور توزيع زنگان'
$ws.Range("B735").Value = 10260057369
$ws.Range("C735").Value = 'توليدي و خدمات صنايع نسوز توكا'
$ws.Range("B736").Value = 10220079472
$ws.Range("C736").Value = 'پتروشيمي اروميه'
$ws.Range("B737").Value = 10220022017
$ws.Range("C737").Value = 'سيمان سفيد اروميه'
$ws.Range("B738").Value = 10103024860
$ws.Range("C738").Value = 'ليزينگ اقتصاد نوين'
$ws.Range("B739").Value = 10380270277
$ws.Range("C739").Value = 'احياء صنايع خراسان'
$ws.Range("B740").Value = 10260326724
$ws.Range("C740").Value = 'م .صنايع و معادن احياء سپاهان'
$ws.Range("B741").Value = 10101324421
$ws.Range("C741").Value = 'سرمايه گذاري صنايع ايران'
$ws.Range("B742").Value = 10100304130
$ws.Range("C742").Value = 'ويتانا'
$ws.Range("B743").Value = 10103972003
$ws.Range("C743").Value = 'واسپاري ملت'
$ws.Range("B744").Value = 10102617399
$ws.Range("C744").Value = 'شركت ليزينگ آريا دانا'
$ws.Range("B745").Value = 14004810068
$ws.Range("C745").Value = 'گ.مديريت ارزش سرمايه ص ب كشوري'
$ws.Range("B746").Value = 10320821816
$ws.Range("C746").Value = 'بانك مهر اقتصاد'
$ws.Range("B747").Value = 10102681950
$ws.Range("C747").Value = 'سرمايه گذاري وثوق امين'
$ws.Range("B748").Value = 10860246171
$ws.Range("C748").Value = 'بانك سينا'
$ws.Range("B749").Value = 10102529006
$ws.Range("C749").Value = 'بيمه سينا'
$ws.Range("B750").Value = 10861638925
$ws.Range("C750").Value = 'فرانسوز يزد'
$ws.Range("B751").Value = 10100440525
$ws.Range("C751").Value = 'ايران‌ياساتايرورابر'
$ws.Range("B752").Value = 10101016508
$ws.Range("C752").Value = 'مجتمع صنايع لاستيك يزد'
$ws.Range("B753").Value = 10102967236
$ws.Range("C753").Value = 'پتروشيمي مرجان'
$ws.Range("B754").Value = 10000000005
$ws.Range("C754").Value = 'عمران و مسكن سازان شمالغرب'
$ws.Range("B755").Value = 10000000006
$ws.Range("C755").Value = 'سرمايه گذاري حافظ اعتماد'
$ws.Range("B756").Value = 10102801066
$ws.Range("C756").Value = 'سرمايه گذاري ارشك'
$ws.Range("B757").Value = 10000000007
$ws.Range("C757").Value = 'سرمايه گذاري سيراف'
$ws.Range("B758").Value = 10102801961
$ws.Range("C758").Value = 'سرمايه گذاري سليم'
$ws.Range("B759").Value = 10102773137
$ws.Range("C759").Value = 'سرمايه گذاري زعيم'
$ws.Range("B760").Value = 10000000008
$ws.Range("C760").Value = 'سرمايه گذاري اعتصام'
$ws.Range("B761").Value = 10000000009
$ws.Range("C761").Value = 'سرمايه گذاري مفتاح'
$ws.Range("B762").Value = 10320856814
$ws.Range("C762").Value = 'مولدنيروگاهي تجارت فارس'
$ws.Range("B763").Value = 10102399345
$ws.Range("C763").Value = 'گسترش صنعت علوم زيستي'
$ws.Range("B764").Value = 10320814852
$ws.Range("C764").Value = 'توليد برق پرند مپنا'
$ws.Range("B765").Value = 10102169938
$ws.Range("C765").Value = 'سرمايه گذاري اعتضاد غدير'
$ws.Range("B766").Value = 10000000010
$ws.Range("C766").Value = 'شركت اعتباري ثامن'
$ws.Range("B767").Value = 10000000011
$ws.Range("C767").Value = 'خدمات مديريت صندوق بازنشستگي'
$ws.Range("B768").Value = 10102855345
$ws.Range("C768").Value = 'خدمات هوايي سامان'
$ws.Range("B769").Value = 10103754436
$ws.Range("C769").Value = 'فولاد زرند ايرانيان'
$ws.Range("B770").Value = 10104088225
$ws.Range("C770").Value = 'فولاد سيرجان ايرانيان'
$ws.Range("B771").Value = 10000000012
$ws.Range("C771").Value = 'تامين سرمايه كيميا'
$ws.Range("B772").Value = 10460086306
$ws.Range("C772").Value = 'صنعت روي زنگان'
$ws.Range("B773").Value = 10100047773
$ws.Range("C773").Value = 'صنعتي زر ماكارون'
$ws.Range("B774").Value = 10101335387
$ws.Range("C774").Value = 'مرغ اجداد زربال'
$ws.Range("B775").Value = 10000000014
$ws.Range("C775").Value = 'توسعه گردشگري ايران'
$ws.Range("B776").Value = 10000000013
$ws.Range("C776").Value = 'حمل و نقل ايران و روسيه'
$ws.Range("B777").Value = 10000000015
$ws.Range("C777").Value = 'مهندسي جوش ايران'
$ws.Range("B778").Value = 10000000016
$ws.Range("C778").Value = 'خدمات غير صنعتي گاز ايران'
$ws.Range("B779").Value = 10000000017
$ws.Range("C779").Value = 'پايانه ها و مخازن پتروشيمي'
$ws.Range("B780").Value = 10000000018
$ws.Range("C780").Value = 'كالاي پتروشيمي'
$ws.Range("B781").Value = 10861351572
$ws.Range("C781").Value = 'توليد نيروي برق سهند'
$ws.Range("B782").Value = 10101336772
$ws.Range("C782").Value = 'بهره برداري نيروگاه دز'
$ws.Range("B783").Value = 14008126329
$ws.Range("C783").Value = 'فن آوا كارت'
$ws.Range("B784").Value = 10103978066
$ws.Range("C784").Value = 'پتروشيمي سلمان فارسي'
$ws.Range("B785").Value = 10100627111
$ws.Range("C785").Value = 'نوسازي صنايع ايران'
$ws.Range("B786").Value = 10100655069
$ws.Range("C786").Value = 'ساختماني عمران تكلار'
$ws.Range("B787").Value = 10480037917
$ws.Range("C787").Value = 'زغالسنگ البرز شرقي'
$ws.Range("B788").Value = 10840007960
$ws.Range("C788").Value = 'سنگ آهن مركزي'
$ws.Range("B789").Value = 14005155456
$ws.Range("C789").Value = 'نيروگاه زاگرس كوثر'
$ws.Range("B790").Value = 10680047080
$ws.Range("C790").Value = 'شير و گوشت زاگرس شهركرد'
$ws.Range("B791").Value = 10320503079
$ws.Range("C791").Value = 'بانك ايران زمين'
$ws.Range("B792").Value = 10100370594
$ws.Range("C792").Value = 'زامياد'
$ws.Range("B793").Value = 10260482101
$ws.Range("C793").Value = 'سرمايه‌گذاري مسكن زاينده رود'
$ws.Range("B794").Value = 10460101525
$ws.Range("C794").Value = 'صنايع كشاورزي وكود زنجان'
$ws.Range("B795").Value = 10260085303
$ws.Range("C795").Value = 'ذوب آهن اصفهان'
$ws.Range("B796").Value = 10240011565
$ws.Range("C796").Value = 'ملي كشت و صنعت و دامپروري پارس'
$ws.Range("B797").Value = 10184001687
$ws.Range("C797").Value = 'گروه ص. پژوهشي فرهيختگان زرنام'
$ws.Range("B798").Value = 10320311482
$ws.Range("C798").Value = 'صرافی فراز اعتماد'
$ws.Range("B799").Value = 10320718470
$ws.Range("C799").Value = 'شبکه الکترونیکی پرداخت کارت شاپرک'
$ws.Range("B800").Value = 10103714248
$ws.Range("C800").Value = 'ساماندهی مطالبات معوق'
$ws.Range("B801").Value = 14008341286
$ws.Range("C801").Value = 'مدیریت بازار متشکل معاملات ارزی'
$ws.Range("B802").Value = 14005790156
$ws.Range("C802").Value = 'سامانه های کاربردی کلان همگام'
$ws.Range("B803").Value = 10102801942
$ws.Range("C803").Value = 'گروه مالی ملت'
$ws.Range("B804").Value = 10101165910
$ws.Range("C804").Value = 'بهساز مشارکت های ملت'
$ws.Range("B805").Value = 14004003162
$ws.Range("C805").Value = 'گروه فن آوران بهسازان فردا'
$ws.Range("B806").Value = 10861406488
$ws.Range("C806").Value = 'تولیدی گرانول قزوین'
$ws.Range("B807").Value = 10102363667
$ws.Range("C807").Value = 'کوبل دارو'
$ws.Range("B808").Value = 14006408218
$ws.Range("C808").Value = 'گروه مدیریت سرمایه لیان'
$ws.Range("B809").Value = 14007070869
$ws.Range("C809").Value = 'آرین الوند پارس'
$ws.Range("B810").Value = 10000000054
$ws.Range("C810").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی صبا گستر نفت و گاز تامین'
$ws.Range("B811").Value = 10260181270
$ws.Range("C811").Value = 'ایثار فجر کاشان'
$ws.Range("B812").Value = 10000000019
$ws.Range("C812").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی آواری زاگرس'
$ws.Range("B813").Value = 10102773194
$ws.Range("C813").Value = 'سرمایه گذاری امین اعتماد'
$ws.Range("B814").Value = 10103492170
$ws.Range("C814").Value = 'گسترش الکترونیک تدبیر ایران'
$ws.Range("B815").Value = 10000000063
$ws.Range("C815").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی تراز ویستا'
$ws.Range("B816").Value = 10103891706
$ws.Range("C816").Value = 'توسعه سرمایه گذاری دریا ساحل ایرانیان'
$ws.Range("B817").Value = 10102987556
$ws.Range("C817").Value = 'پیشتازان تجارت ویستا'
$ws.Range("B818").Value = 10103891690
$ws.Range("C818").Value = 'توسعه تجارت مجازی سارینا'
$ws.Range("B819").Value = 10102911447
$ws.Range("C819").Value = 'نو آوری ستاره پارسیان'
$ws.Range("B820").Value = 14008302538
$ws.Range("C820").Value = 'پیشگامان رشد و نوآوری'
$ws.Range("B821").Value = 10103137408
$ws.Range("C821").Value = 'سرمایه گذاری آتیه صبا'
$ws.Range("B822").Value = 14010024250
$ws.Range("C822").Value = 'توسعه کارآفرینی ماکان'
$ws.Range("B823").Value = 14006953961
$ws.Range("C823").Value = 'دنیای رشد و نوآوری'
$ws.Range("B824").Value = 14009224728
$ws.Range("C824").Value = 'رادیس'
$ws.Range("B825").Value = 10101885990
$ws.Range("C825").Value = 'صنعتی آراسته معدن'
$ws.Range("B826").Value = 10861503629
$ws.Range("C826").Value = 'کوشش آذین قشم'
$ws.Range("B827").Value = 10861394039
$ws.Range("C827").Value = 'سرمایه گذاری گلومینکو قشم'
$ws.Range("B828").Value = 14001907528
$ws.Range("C828").Value = 'بیمه مرکزی جمهوری اسلامی ایران'
$ws.Range("B829").Value = 14004046372
$ws.Range("C829").Value = 'گروه مالی پارسیان'
$ws.Range("B830").Value = 10103652364
$ws.Range("C830").Value = 'صنایع آلوم رول نوین'
$ws.Range("B831").Value = 10780110885
$ws.Range("C831").Value = 'روان گداز پردیس'
$ws.Range("B832").Value = 10103529793
$ws.Range("C832").Value = 'مدبر تجارت آریا'
$ws.Range("B833").Value = 10103529800
$ws.Range("C833").Value = 'پرتو کالا پردیس'
$ws.Range("B834").Value = 10101195187
$ws.Range("C834").Value = 'امین آر'
$ws.Range("B835").Value = 10380430077
$ws.Range("C835").Value = 'سرمایه گذاری سهام عدالت استان خراسان رضوی'
$ws.Range("B836").Value = 10530320996
$ws.Range("C836").Value = 'سرمایه گذاری استان فارس'
$ws.Range("B837").Value = 10860903365
$ws.Range("C837").Value = 'سرمایه گذاری استان خوزستان'
$ws.Range("B838").Value = 10260490269
$ws.Range("C838").Value = 'سرمایه گذاری استان اصفهان'
$ws.Range("B839").Value = 10102158022
$ws.Range("C839").Value = 'سازمان توسعه و نوسازی معادن و صنایع معدنی ایران'
$ws.Range("B840").Value = 14008307703
$ws.Range("C840").Value = 'سرو سودمند مدبران'
$ws.Range("B841").Value = 14006249920
$ws.Range("C841").Value = 'تلاش انگیزه دارای آریا'
$ws.Range("B842").Value = 14011339503
$ws.Range("C842").Value = 'امین سلامت بهاران'
$ws.Range("B843").Value = 10320866327
$ws.Range("C843").Value = 'صندوق سرمایه گذاری اندوخته پایدار سپهر'
$ws.Range("B844").Value = 14006852814
$ws.Range("C844").Value = 'سازمان اقتصادی کوثر'
$ws.Range("B845").Value = 14006012357
$ws.Range("C845").Value = 'اختصاصی بازارگردانی گروه دی'
$ws.Range("B846").Value = 14008473366
$ws.Range("C846").Value = 'اختصاصی بازارگردانی پاداش پشتیبان پارس'
$ws.Range("B847").Value = 10000000059
$ws.Range("C847").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی توسعه بازار تمدن'
$ws.Range("B848").Value = 14003959426
$ws.Range("C848").Value = 'رسا سامان آریا'
$ws.Range("B849").Value = 10320600788
$ws.Range("C849").Value = 'صندوق سرمایه گذاری بانک گردشگری'
$ws.Range("B850").Value = 14005496342
$ws.Range("C850").Value = 'صندوق سرمايه گذاري اندوخته توسعه صادرات آرماني'
$ws.Range("B851").Value = 10000000099
$ws.Range("C851").Value = 'صندوق سرمايه گذاري امين يكم فردا'
$ws.Range("B852").Value = 14004128084
$ws.Range("C852").Value = 'پارس تامین مجد'
$ws.Range("B853").Value = 10103679112
$ws.Range("C853").Value = 'گسترش فناوری های نوین'
$ws.Range("B854").Value = 10100589211
$ws.Range("C854").Value = 'رفاه وتامین اتیه امید'
$ws.Range("B855").Value = 10320752394
$ws.Range("C855").Value = 'صندوق سرمایه گذاری امین آشنا ایرانیان'
$ws.Range("B856").Value = 14004183170
$ws.Range("C856").Value = 'حکمت آشنا ایرانیان'
$ws.Range("B857").Value = 10100258391
$ws.Range("C857").Value = 'ملی نفتکش ایران'
$ws.Range("B858").Value = 10100906047
$ws.Range("C858").Value = 'حمل واردات ایران'
$ws.Range("B859").Value = 10100976112
$ws.Range("C859").Value = 'حمل دریایی ایران'
$ws.Range("B860").Value = 10103372005
$ws.Range("C860").Value = 'مفید شیشه'
$ws.Range("B861").Value = 14007120347
$ws.Range("C861").Value = 'اختصاصی بازارگردانی سهم آشنا یکم'
$ws.Range("B862").Value = 14009518268
$ws.Range("C862").Value = 'نام آوران تجارت اردستان'
$ws.Range("B863").Value = 10630025764
$ws.Range("C863").Value = 'صندوق بازنشستگی شرکت ملی صنایع مس ایران'
$ws.Range("B864").Value = 10101915348
$ws.Range("C864").Value = 'کارخانجات نورد لوله یاران'
$ws.Range("B865").Value = 10320469990
$ws.Range("C865").Value = 'پاتین بین الملل پردیس'
$ws.Range("B866").Value = 14007123145
$ws.Range("C866").Value = 'تامین آتیه سرزمین ایرانیان'
$ws.Range("B867").Value = 10103258569
$ws.Range("C867").Value = 'بين المللي سرمايه گذاري ايرانيان'
$ws.Range("B868").Value = 10380450359
$ws.Range("C868").Value = 'سرمایه گذاری آینده نگر شرق'
$ws.Range("B869").Value = 14006336023
$ws.Range("C869").Value = 'عمران پی ژیوار'
$ws.Range("B870").Value = 10320400445
$ws.Range("C870").Value = 'طلایه داران تجارت کاسپین'
$ws.Range("B871").Value = 14008157950
$ws.Range("C871").Value = 'نیکان سرشت سرزمین ایرانیان'
$ws.Range("B872").Value = 10861529184
$ws.Range("C872").Value = 'آی اف اس کیش'
$ws.Range("B873").Value = 10000000029
$ws.Range("C873").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی افتخار حافظ'
$ws.Range("B874").Value = 10102803281
$ws.Range("C874").Value = 'سرمایه گذاری دلیران پارس'
$ws.Range("B875").Value = 10320650723
$ws.Range("C875").Value = 'پیشگامان توسعه شهری ایرانیان'
$ws.Range("B876").Value = 10380305808
$ws.Range("C876").Value = 'آساگستران خراسان'
$ws.Range("B877").Value = 14009263200
$ws.Range("C877").Value = 'سرمایه گذاری کشاورزی آرتین هیواد'
$ws.Range("B878").Value = 10000000060
$ws.Range("C878").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی توسعه تاک دانا'
$ws.Range("B879").Value = 10320813657
$ws.Range("C879").Value = 'صندوق سرمایه گذاری مشترک آسمان یکم'
$ws.Range("B880").Value = 10102694700
$ws.Range("C880").Value = 'سرمایه گذاری گروه صنعتی رازی'
$ws.Range("B881").Value = 14005108893
$ws.Range("C881").Value = 'اختصاصی بازار گردان صبا نیک'
$ws.Range("B882").Value = 10000000036
$ws.Range("C882").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی توسعه فولاد مبارکه'
$ws.Range("B883").Value = 14000074621
$ws.Range("C883").Value = 'شرکت پویا فراز کیش'
$ws.Range("B884").Value = 10861530138
$ws.Range("C884").Value = 'سرمایه گذاری ای اف جی اج کیش'
$ws.Range("B885").Value = 10103037570
$ws.Range("C885").Value = 'مؤسسه صندوق بازنشستگی وظیفه از کار افتادگی و پس انداز کارکنان بانکهای ملی و ادغام شده'
$ws.Range("B886").Value = 10101997560
$ws.Range("C886").Value = 'صبا میهن'
$ws.Range("B887").Value = 14004107956
$ws.Range("C887").Value = 'مشترک ارزش کاوان آینده'
$ws.Range("B888").Value = 10380520955
$ws.Range("C888").Value = 'توسعه توس بان امین'
$ws.Range("B889").Value = 14006243544
$ws.Range("C889").Value = 'گسترش سرمایه هوشمند'
$ws.Range("B890").Value = 10103761233
$ws.Range("C890").Value = 'گروه مالی ملل'
$ws.Range("B891").Value = 14004407347
$ws.Range("C891").Value = 'بازارگردانی نوین پیشرو'
$ws.Range("B892").Value = 10380647203
$ws.Range("C892").Value = 'افق رهباد خاوران توس'
$ws.Range("B893").Value = 14007411551
$ws.Range("C893").Value = 'تامین انرژی سپنتا توس'
$ws.Range("B894").Value = 10320352137
$ws.Range("C894").Value = 'گروه توسعه اقتصاد ملل'
$ws.Range("B895").Value = 10320774211
$ws.Range("C895").Value = 'تجارت نصر البرز'
$ws.Range("B896").Value = 10380647222
$ws.Range("C896").Value = 'خاوران جم گسترش تابران'
$ws.Range("B897").Value = 10862064732
$ws.Range("C897").Value = 'بانک سپه'
$ws.Range("B898").Value = 10101364297
$ws.Range("C898").Value = 'صنعتی و بازرگانی غدیر'
$ws.Range("B899").Value = 10103837742
$ws.Range("C899").Value = 'سرمایه گذاری آذر'
$ws.Range("B900").Value = 10000000043
$ws.Range("C900").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی لاجورد دماوند'
$ws.Range("B901").Value = 10102898876
$ws.Range("C901").Value = 'سرمایه گذاری مهرگان سرمایه پارس'
$ws.Range("B902").Value = 10000000024
$ws.Range("C902").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی آتیه باران'
$ws.Range("B903").Value = 14008247684
$ws.Range("C903").Value = 'فناوری اطلاعات و ارتباطات راهبرد'
$ws.Range("B904").Value = 14007707492
$ws.Range("C904").Value = 'رسا سازه پی کاوان مهام'
$ws.Range("B905").Value = 10320692265
$ws.Range("C905").Value = 'آرتا ایده نفیس'
$ws.Range("B906").Value = 10320673140
$ws.Range("C906").Value = 'سرمایه گذاری صنعتی معدنی آریا فاتح خاورمیانه'
$ws.Range("B907").Value = 10320692094
$ws.Range("C907").Value = 'سرمایه گذاری آتیه مداران'
$ws.Range("B908").Value = 14004068870
$ws.Range("C908").Value = 'خدمات مدیریت اندیشه آتیه مداران'
$ws.Range("B909").Value = 14008556974
$ws.Range("C909").Value = 'فناوری ارتباطات و اطلاعات ایده دیجیتال هوشمند'
$ws.Range("B910").Value = 10103460290
$ws.Range("C910").Value = 'سرمایه گذاری مهرگان تامین پارس'
$ws.Range("B911").Value = 10101442691
$ws.Range("C911").Value = 'کارگزاری بانک صنعت ومعدن'
$ws.Range("B912").Value = 10320307426
$ws.Range("C912").Value = 'سرمایه گذاری تامین آتیه مسکن'
$ws.Range("B913").Value = 14005253222
$ws.Range("C913").Value = 'صنعتی ومعدنی ایران'
$ws.Range("B914").Value = 10200450670
$ws.Range("C914").Value = 'گروه صنعتی اشتالز فولاد خاورمیانه'
$ws.Range("B915").Value = 14008076836
$ws.Range("C915").Value = 'تدبیر فردای نیک'
$ws.Range("B916").Value = 10100937493
$ws.Range("C916").Value = 'مدیریت طرح و توسعه آینده پویا'
$ws.Range("B917").Value = 10100505148
$ws.Range("C917").Value = 'مؤسسه بهارستان آسایش'
$ws.Range("B918").Value = 10101814661
$ws.Range("C918").Value = 'آرتا مبین سحر'
$ws.Range("B919").Value = 10320884761
$ws.Range("C919").Value = 'محب گردشگری سلامت ایرانیان'
$ws.Range("B920").Value = 10102173008
$ws.Range("C920").Value = 'گروه تولیدی رنان طب'
$ws.Range("B921").Value = 10102593722
$ws.Range("C921").Value = 'گسترش فناوری عمران زیست'
$ws.Range("B922").Value = 10104060856
$ws.Range("C922").Value = 'محب مهر سلامت پارس'
$ws.Range("B923").Value = 10000000034
$ws.Range("C923").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی توسعه سهام نیکی'
$ws.Range("B924").Value = 10000000033
$ws.Range("C924").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی تصمیم ساز'
$ws.Range("B925").Value = 10320895923
$ws.Range("C925").Value = 'سرزمین پهناور مهر'
$ws.Range("B926").Value = 14005828941
$ws.Range("C926").Value = 'مؤسسه دانش بنیان برکت'
$ws.Range("B927").Value = 14006793310
$ws.Range("C927").Value = 'اختصاصی بازارگردانی تدبیرگران فردا'
$ws.Range("B928").Value = 14005786478
$ws.Range("C928").Value = 'مدیریت سرمایه ارزش آفرین دانا'
$ws.Range("B929").Value = 14005965618
$ws.Range("C929").Value = 'مدیریت سرمایه آسای دانا'
$ws.Range("B930").Value = 14005925472
$ws.Range("C930").Value = 'توسعه سرمایه پارمیس'
$ws.Range("B931").Value = 14007027064
$ws.Range("C931").Value = 'باز آفرین ابهر نیکو'
$ws.Range("B932").Value = 14007137903
$ws.Range("C932").Value = 'توسعه کسب و کار باتیس'
$ws.Range("B933").Value = 14006082952
$ws.Range("C933").Value = 'راه سازان تلاش معدن'
$ws.Range("B934").Value = 10320780362
$ws.Range("C934").Value = 'بنیاد خیریه تات'
$ws.Range("B935").Value = 10320858723
$ws.Range("C935").Value = 'زینت تجارت آریا'
$ws.Range("B936").Value = 10320661767
$ws.Range("C936").Value = 'پارس سرمایه تابا'
$ws.Range("B937").Value = 10320169213
$ws.Range("C937").Value = 'تولیدی و صنعتی مهرآوران آتیه البرز'
$ws.Range("B938").Value = 10000000055
$ws.Range("C938").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی امیدلوتوس پارسیان'
$ws.Range("B939").Value = 14005933591
$ws.Range("C939").Value = 'گروه مالی و اقتصادی آینده'
$ws.Range("B940").Value = 10980226256
$ws.Range("C940").Value = 'تجارت کالای دنا کیش'
$ws.Range("B941").Value = 10320503955
$ws.Range("C941").Value = 'صندوق سرمایه گذاری یکم کارگزاری بانک کشاورزی'
$ws.Range("B942").Value = 10000000026
$ws.Range("C942").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی آرمان اندیش'
$ws.Range("B943").Value = 10420196707
$ws.Range("C943").Value = 'فنی و مهندسی جنوب تاسیسات'
$ws.Range("B944").Value = 10860411138
$ws.Range("C944").Value = 'مهندسی بازرگانی پیام قشم'
$ws.Range("B945").Value = 10101759580
$ws.Range("C945").Value = 'خدمات بازرگانی پیمان امیر'
$ws.Range("B946").Value = 10320832748
$ws.Range("C946").Value = 'گروه مالی شهر'
$ws.Range("B947").Value = 10760094603
$ws.Range("C947").Value = 'گسترش کشاورزی و دامپروری فردوس پارس'
$ws.Range("B948").Value = 10100171920
$ws.Range("C948").Value = 'بنیاد مستضعفان انقلاب اسلامی'
$ws.Range("B949").Value = 10000000040
$ws.Range("C949").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی سینا بهگزین'
$ws.Range("B950").Value = 10100625070
$ws.Range("C950").Value = 'مؤسسه صندوق بیمه اجتماعی روستاییان و عشایر'
$ws.Range("B951").Value = 14003674869
$ws.Range("C951").Value = 'خردمندان صابر عصر'
$ws.Range("B952").Value = 14009108895
$ws.Range("C952").Value = 'صندوق سرمايه گذاري واسطه گري مالي يكم'
$ws.Range("B953").Value = 10000000030
$ws.Range("C953").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی اکسیرسودا'
$ws.Range("B954").Value = 14008036097
$ws.Range("C954").Value = 'سرمایه گذاری تجاری سامانه های رایانه ای هوپاد هونامیک'
$ws.Range("B955").Value = 10100053930
$ws.Range("C955").Value = 'خدمات مدیریت صبا تامین'
$ws.Range("B956").Value = 14004684611
$ws.Range("C956").Value = 'با درآمد ثابت کاردان'
$ws.Range("B957").Value = 10100351811
$ws.Range("C957").Value = 'گروه مالی و اقتصادی دی'
$ws.Range("B958").Value = 10320539175
$ws.Range("C958").Value = 'خدمات مالی حسابداری دی ایرانیان'
$ws.Range("B959").Value = 10100567400
$ws.Range("C959").Value = 'اندوخته شاهد'
$ws.Range("B960").Value = 14009263160
$ws.Range("C960").Value = 'سرمایه گذاری اقتصادی هامرز راتین'
$ws.Range("B961").Value = 10861677542
$ws.Range("C961").Value = 'بانک ملی'
$ws.Range("B962").Value = 10000000073
$ws.Range("C962").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی توسعه ملل'
$ws.Range("B963").Value = 10103858742
$ws.Range("C963").Value = 'سهامی بیمه ایران'
$ws.Range("B964").Value = 10102279261
$ws.Range("C964").Value = 'گروه مالی سپه'
$ws.Range("B965").Value = 10102803034
$ws.Range("C965").Value = 'سرمايه گذاری يمين'
$ws.Range("B966").Value = 14003699653
$ws.Range("C966").Value = 'صندوق بازنشستگی کارکنان صدا و سیما'
$ws.Range("B967").Value = 14004081532
$ws.Range("C967").Value = 'ارزش آفرینان صبا'
$ws.Range("B968").Value = 10103297319
$ws.Range("C968").Value = 'سرمايه گذاري نگين خاتم ايرانيان'
$ws.Range("B969").Value = 10380179815
$ws.Range("C969").Value = 'مکین'
$ws.Range("B970").Value = 14003902009
$ws.Range("C970").Value = 'صندوق سرمایه گذاری مشترک بانک خاورمیانه'
$ws.Range("B971").Value = 10200258715
$ws.Range("C971").Value = 'سرمایه گذاری استان آذربایجان شرقی'
$ws.Range("B972").Value = 10760345687
$ws.Range("C972").Value = 'سرمایه گذاری استان مازندران'
$ws.Range("B973").Value = 10630135824
$ws.Range("C973").Value = 'سرمایه گذاری استان کرمان'
$ws.Range("B974").Value = 10220124165
$ws.Range("C974").Value = 'سرمایه گذاری استان آذربایجان غربی'
$ws.Range("B975").Value = 10860405350
$ws.Range("C975").Value = 'سرمایه گذاری استان سیستان و بلوچستان'
$ws.Range("B976").Value = 10720243977
$ws.Range("C976").Value = 'سرمایه گذاری استان گیلان'
$ws.Range("B977").Value = 10700141201
$ws.Range("C977").Value = 'سرمایه گذاری استان گلستان'
$ws.Range("B978").Value = 10320528091
$ws.Range("C978").Value = 'آتیه سازان دی'
$ws.Range("B979").Value = 10320354635
$ws.Range("C979").Value = 'سرمایه گذاری مدبران اقتصاد'
$ws.Range("B980").Value = 10101380313
$ws.Range("C980").Value = 'سرمایه گذاری فرهنگیان'
$ws.Range("B981").Value = 10103349751
$ws.Range("C981").Value = 'گروه سرمایه گذاری فولاد گستر کوثر'
$ws.Range("B982").Value = 10104052682
$ws.Range("C982").Value = 'توسعه سرمایه گذاری دریک'
$ws.Range("B983").Value = 14004356528
$ws.Range("C983").Value = 'مجتمع آهن و فولاد قایم سبلان'
$ws.Range("B984").Value = 10260434191
$ws.Range("C984").Value = 'آزمون احیاء سپاهان'
$ws.Range("B985").Value = 10260334425
$ws.Range("C985").Value = 'مهندسی قایم سپاهان'
$ws.Range("B986").Value = 10102399005
$ws.Range("C986").Value = 'گروه مهدتاژ'
$ws.Range("B987").Value = 10103663870
$ws.Range("C987").Value = 'گروه توسعه ساختمان تدبیر'
$ws.Range("B988").Value = 10103871690
$ws.Range("C988").Value = 'گروه توسعه اقتصادی تدبیر'
$ws.Range("B989").Value = 10320823442
$ws.Range("C989").Value = 'صندوق سرمایه گذاری لوتوس پارسیان'
$ws.Range("B990").Value = 14006847020
$ws.Range("C990").Value = 'با درآمد ثابت کمند'
$ws.Range("B991").Value = 14003772673
$ws.Range("C991").Value = 'صندوق سرمایه گذاری مشترک سپهر تدبیرگران'
$ws.Range("B992").Value = 10104055685
$ws.Range("C992").Value = 'مادر تخصصی گسترش صنایع غذایی سینا'
$ws.Range("B993").Value = 14006123572
$ws.Range("C993").Value = 'توسعه بازاریابی و فروش برق صبا'
$ws.Range("B994").Value = 10000000038
$ws.Range("C994").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی خلیج فارس'
$ws.Range("B995").Value = 10102986880
$ws.Range("C995").Value = 'سرمایه گذاری زرین پرشیا'
$ws.Range("B996").Value = 10000000028
$ws.Range("C996").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی آرمان انصار'
$ws.Range("B997").Value = 10000000050
$ws.Range("C997").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی مهرگان'
$ws.Range("B998").Value = 10100845276
$ws.Range("C998").Value = 'تعاونی مصرف کارکنان ستاد ارتش جمهوری اسلامی ایران'
$ws.Range("B999").Value = 10101727833
$ws.Range("C999").Value = 'خدمات هواپیمایی کاسپین'
$ws.Range("B1000").Value = 10840070153
$ws.Range("C1000").Value = 'تعاونی خدماتی پیشگامان کویر یزد'
$ws.Range("B1001").Value = 10840086913
$ws.Range("C1001").Value = 'پیشگامان کویر آسیا'
$ws.Range("B1002").Value = 10000000046
$ws.Range("C1002").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی بهمن گست'
$ws.Range("B1003").Value = 10102685249
$ws.Range("C1003").Value = 'سرمایه گذاری آوین'
$ws.Range("B1004").Value = 14003725635
$ws.Range("C1004").Value = 'صندوق سرمایه گذاری آسمان آرمانی سهام'
$ws.Range("B1005").Value = 14003767221
$ws.Range("C1005").Value = 'صندوق سرمایه گذاری توسعه اندوخته آینده'
$ws.Range("B1006").Value = 10102727560
$ws.Range("C1006").Value = 'فن پردازان بهمن'
$ws.Range("B1007").Value = 10102803015
$ws.Range("C1007").Value = 'سرمایه گذاری اندیشه فردا'
$ws.Range("B1008").Value = 10102802996
$ws.Range("C1008").Value = 'سرمایه گذاری توسعه انرژی هور'
$ws.Range("B1009").Value = 10102803205
$ws.Range("C1009").Value = 'سرمایه گذاری مانا نوین'
$ws.Range("B1010").Value = 10101997272
$ws.Range("C1010").Value = 'بازرگانی عصر بهمن'
$ws.Range("B1011").Value = 10000000041
$ws.Range("C1011").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی گروه توسعه بهشهر'
$ws.Range("B1012").Value = 10000000075
$ws.Range("C1012").Value = 'بيمه مركزي ايران-صندوق بازنشستگي'
$ws.Range("B1013").Value = 10000000092
$ws.Range("C1013").Value = 'صندوق بازنشستگي کارکنان بانک مرکزی جمهوری اسلامي ايران'
$ws.Range("B1014").Value = 10102960400
$ws.Range("C1014").Value = 'توسعه مدیریت پارس حافظ'
$ws.Range("B1015").Value = 10103705170
$ws.Range("C1015").Value = 'نظم آوران شایسته'
$ws.Range("B1016").Value = 10320839039
$ws.Range("C1016").Value = 'نسیم تجارت فردا'
$ws.Range("B1017").Value = 10320836357
$ws.Range("C1017").Value = 'توسعه تجارت غرب ایرانیان'
$ws.Range("B1018").Value = 10320876676
$ws.Range("C1018").Value = 'ارمغان تجارت پایدار'
$ws.Range("B1019").Value = 10320635617
$ws.Range("C1019").Value = 'تجارت و توسعه مهرآفرینان ونداد'
$ws.Range("B1020").Value = 10320635904
$ws.Range("C1020").Value = 'نظم آوران پویای ایرانیان'
$ws.Range("B1021").Value = 14004210872
$ws.Range("C1021").Value = 'نیک اندیشان سروش فجر'
$ws.Range("B1022").Value = 14004156637
$ws.Range("C1022").Value = 'دانش گستران آتی ساز پویا'
$ws.Range("B1023").Value = 14008079936
$ws.Range("C1023").Value = 'همیاری سرمایه انسانی پاسارگاد'
$ws.Range("B1024").Value = 10320634937
$ws.Range("C1024").Value = 'تجارت و توسعه مهرآفرینان سرآمد'
$ws.Range("B1025").Value = 10100385648
$ws.Range("C1025").Value = 'بنیاد فرهنگی مصلی نژاد'
$ws.Range("B1026").Value = 10320635047
$ws.Range("C1026").Value = 'اندیشه سازان بسامان ونداد'
$ws.Range("B1027").Value = 10630107136
$ws.Range("C1027").Value = 'معدنکاری اولنگ'
$ws.Range("B1028").Value = 10380474989
$ws.Range("C1028").Value = 'کیا آسا تجارت توس'
$ws.Range("B1029").Value = 14004975980
$ws.Range("C1029").Value = 'مشترک افق کارگزاری بانک خاورمیانه'
$ws.Range("B1030").Value = 10102417191
$ws.Range("C1030").Value = 'سام گروه'
$ws.Range("B1031").Value = 14009224859
$ws.Range("C1031").Value = 'گروه راما'
$ws.Range("B1032").Value = 10101730772
$ws.Range("C1032").Value = 'کارگزاری سی ولکس'
$ws.Range("B1033").Value = 10102558446
$ws.Range("C1033").Value = 'کارگزاری ستاره جنوب'
$ws.Range("B1034").Value = 10101476080
$ws.Range("C1034").Value = 'کارگزاری بانک کشاورزی'
$ws.Range("B1035").Value = 10760335630
$ws.Range("C1035").Value = 'کارگزاری بانک آینده'
$ws.Range("B1036").Value = 10102702966
$ws.Range("C1036").Value = 'کارگزاری صبا جهاد'
$ws.Range("B1037").Value = 10102634872
$ws.Range("C1037").Value = 'کارگزاری توسعه کشاورزی'
$ws.Range("B1038").Value = 10102002047
$ws.Range("C1038").Value = 'کارگزاری نهایت نگر'
$ws.Range("B1039").Value = 10860222183
$ws.Range("C1039").Value = 'کارگزاری مبین سرمایه'
$ws.Range("B1040").Value = 10102775538
$ws.Range("C1040").Value = 'کارگزاری سینا'
$ws.Range("B1041").Value = 14007297814
$ws.Range("C1041").Value = 'گروه خدمات بازار سرمایه مفید'
$ws.Range("B1042").Value = 10101559627
$ws.Range("C1042").Value = 'کارگزاری صبا تامین'
$ws.Range("B1043").Value = 10102040894
$ws.Range("C1043").Value = 'کارگزاری آتی ساز بازار'
$ws.Range("B1044").Value = 10102676190
$ws.Range("C1044").Value = 'پرديس متحد آريا'
$ws.Range("B1045").Value = 10102631183
$ws.Range("C1045").Value = 'كارگزاري كالاي كشاورزي پيمان گستر'
$ws.Range("B1046").Value = 10101553125
$ws.Range("C1046").Value = 'کارگزاری خبرگان سهام'
$ws.Range("B1047").Value = 10104011299
$ws.Range("C1047").Value = 'توانمند سازی بازنشستگان نیروهای مسلح'
$ws.Range("B1048").Value = 10103653560
$ws.Range("C1048").Value = 'سرمایه گذاری تجاری شستان'
$ws.Range("B1049").Value = 10000000093
$ws.Range("C1049").Value = 'صندوق بيمه عمرپرسنل سپاه'
$ws.Range("B1050").Value = 10000000094
$ws.Range("C1050").Value = 'صندوق بيمه عمرپرسنل ناجا'
$ws.Range("B1051").Value = 10000000095
$ws.Range("C1051").Value = 'صندوق بيمه عمرپرسنل آجا'
$ws.Range("B1052").Value = 10101954690
$ws.Range("C1052").Value = 'خدمات بیمه ای سهند مشاور'
$ws.Range("B1053").Value = 10320330972
$ws.Range("C1053").Value = 'سرمایه گذاری الماس حکمت ایرانیان'
$ws.Range("B1054").Value = 10320683865
$ws.Range("C1054").Value = 'آینده سازان فردای کوثر'
$ws.Range("B1055").Value = 10320145605
$ws.Range("C1055").Value = 'مؤسسه بنیاد تعاون وزارت دفاع و پشتیبانی نیروهای مسلح'
$ws.Range("B1056").Value = 10260135725
$ws.Range("C1056").Value = 'مؤسسه فرهنگی خدمات مسافرت هوایی گردشگری و زیارتی ثامن الایمه'
$ws.Range("B1057").Value = 10104000426
$ws.Range("C1057").Value = 'سرمایه آب خاک توسعه'
$ws.Range("B1058").Value = 10380667692
$ws.Range("C1058").Value = 'بازرگانی فرا گستر مهان آرمان'
$ws.Range("B1059").Value = 10102927198
$ws.Range("C1059").Value = 'گروه مالی دانایان'
$ws.Range("B1060").Value = 10860712482
$ws.Range("C1060").Value = 'توسعه بازرگانی هیواد'
$ws.Range("B1061").Value = 14004460426
$ws.Range("C1061").Value = 'هلدینگ سرآمد'
$ws.Range("B1062").Value = 10000000051
$ws.Range("C1062").Value = 'صندوق سرمایه گذاری اختصاصی بازارگردانی دانایان'
$ws.Range("B1063").Value = 10101128860
$ws.Range("C1063").Value = 'پروژه های صنعتی ایران'
$ws.Range("B1064").Value = 10320824078
$ws.Range("C1064").Value = 'پدیده آفرین شفق'
$ws.Range("B1065").Value = 10100687843
$ws.Range("C1065").Value = 'گروه مالی بانک مسکن'
$ws.Range("B1066").Value = 10320782725
$ws.Range("C1066").Value = 'مسکن کارکنان بانک دی'
$ws.Range("B1067").Value = 10320658317
$ws.Range("C1067").Value = 'صندوق سرمایه گذاری ارزش آفرینان دی'
$ws.Range("B1068").Value = 14007562887
$ws.Range("C1068").Value = 'سرمایه گذاری ساختمانی ارزش زمان'
$ws.Range("B1069").Value = 14007982406
$ws.Range("C1069").Value = 'سرمایه گذاری دارویی بهیان پرمون'
$ws.Range("B1070").Value = 14008783355
$ws.Range("C1070").Value = 'سرمایه گذاری معدنی اسپاد تجارت هیوا'
$ws.Range("B1071").Value = 10320637769
$ws.Range("C1071").Value = 'افق نیلی خلیج فارس'
$ws.Range("B1072").Value = 10103679146
$ws.Range("C1072").Value = 'سرمایه گذاری سایه گستر سرمایه'
$ws.Range("B1073").Value = 10200255770
$ws.Range("C1073").Value = 'ویرا سهند تبریز'
$ws.Range("B1074").Value = 14000187509
$ws.Range("C1074").Value = 'دولت جمهوري اسلامي ايران'

# --- Column D (Malekiat_Dolat_Dar_Sherkat) updates ---
$ws.Range("D53").Value = 0.0635
$ws.Range("D59").Value = 0.124
$ws.Range("D86").Value = 0.0992
$ws.Range("D88").Value = 0.0334
$ws.Range("D93").Value = 0.012
$ws.Range("D98").Value = 0
$ws.Range("D102").Value = 0
$ws.Range("D108").Value = 0.063
$ws.Range("D112").Value = 0.05
$ws.Range("D165").Value = 0.0673
$ws.Range("D176").Value = 0.0648
$ws.Range("D201").Value = 0.03
$ws.Range("D255").Value = 0.6596
$ws.Range("D266").Value = 0.0236
$ws.Range("D350").Value = 0.0751
$ws.Range("D406").Value = 0.0738
$ws.Range("D411").Value = 0
$ws.Range("D440").Value = 0.0636
$ws.Range("D540").Value = 0.0891
$ws.Range("D574").Value = 0.2514
$ws.Range("D727").Value = 0.2563
